$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = 1.07
$ws.Range("O2").Value = 1.41
$ws.Range("P2").Value = 2.62
$ws.Range("G3").Value = 2.2
$ws.Range("K3").Value = 2
$ws.Range("U3").Value = 2
$ws.Range("V3").Value = 1.73
$ws.Range("X3").Value = 9.5
$ws.Range("AA3").Value = 21
$ws.Range("AC3").Value = 7.5
$ws.Range("AF3").Value = 67
$ws.Range("AG3").Value = 8.5
$ws.Range("AH3").Value = 15
$ws.Range("AI3").Value = 13
$ws.Range("AX3").Value = 21
$ws.Range("AY3").Value = 34
$ws.Range("G5").Value = 8.25
$ws.Range("H5").Value = 4.55
$ws.Range("I5").Value = 1.32
$ws.Range("J5").Value = 7.4
$ws.Range("K5").Value = 2.4
$ws.Range("L5").Value = 1.8
$ws.Range("N5").Value = 8.25
$ws.Range("P5").Value = 3.85
$ws.Range("Q5").Value = 1.7
$ws.Range("R5").Value = 2.07
$ws.Range("S5").Value = 1.34
$ws.Range("T5").Value = 3
$ws.Range("W5").Value = 22
$ws.Range("X5").Value = 60
$ws.Range("Y5").Value = 27
$ws.Range("AC5").Value = 8.25
$ws.Range("AH5").Value = 6.2
$ws.Range("AN5").Value = 9.25
$ws.Range("AR5").Value = 400
$ws.Range("AT5").Value = 3
$ws.Range("AV5").Value = 90
$ws.Range("AX5").Value = 5.9
$ws.Range("AY5").Value = 17
$ws.Range("AZ5").Value = 16
$ws.Range("G7").Value = 4.5
$ws.Range("H7").Value = 3.4
$ws.Range("I7").Value = 1.7
$ws.Range("J7").Value = 4.85
$ws.Range("K7").Value = 2.15
$ws.Range("L7").Value = 2.27
$ws.Range("N7").Value = 7
$ws.Range("O7").Value = 1.33
$ws.Range("P7").Value = 3.05
$ws.Range("Q7").Value = 1.98
$ws.Range("U7").Value = 1.88
$ws.Range("V7").Value = 1.82
$ws.Range("W7").Value = 11.75
$ws.Range("X7").Value = 26
$ws.Range("Y7").Value = 15
$ws.Range("Z7").Value = 80
$ws.Range("AB7").Value = 50
$ws.Range("AC7").Value = 7
$ws.Range("AD7").Value = 6.8
$ws.Range("AE7").Value = 16.5
$ws.Range("AF7").Value = 80
$ws.Range("AG7").Value = 6.4
$ws.Range("AH7").Value = 7.7
$ws.Range("AJ7").Value = 13.5
$ws.Range("AK7").Value = 14.5
$ws.Range("AM7").Value = 700
$ws.Range("AN7").Value = 6.3
$ws.Range("AO7").Value = 26
$ws.Range("AP7").Value = 32
$ws.Range("AQ7").Value = 150
$ws.Range("AS7").Value = 400
$ws.Range("AU7").Value = 7.5
$ws.Range("AV7").Value = 70
$ws.Range("AW7").Value = 3.55
$ws.Range("AX7").Value = 8.5
$ws.Range("AZ7").Value = 29
